# Add a new "Slovakia" market sheet, based on the existing "Portugal" sheet,
# and populate it with the Slovakia-specific test data.

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Make Portugal the active sheet and select the whole sheet - this mirrors
# the state the source sheet is left in once it stops being the active tab
# after the copy/activation of the new sheet below.
$portugal.Activate()
$portugal.Cells.Select()

# Duplicate the Portugal sheet; the copy is placed immediately after it.
$portugal.Copy([System.Reflection.Missing]::Value, $portugal)

$slovakia = $wb.Worksheets.Item("Portugal (2)")
$slovakia.Name = "Slovakia"

# The source sheet had two rows with an explicit (larger) row height; the
# new sheet's rows should use the default sheet row height instead, so
# auto-fit them back down.
$slovakia.Range("A3:D5").EntireRow.AutoFit()

# The "Input Value" cell for this market loses the bordered box formatting
# that the template cell had.
$slovakia.Range("B4").ClearFormats()

# Market-specific text.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3175/T1827"

# Leave the new sheet active, with B4 selected, matching how the workbook
# was left after adding the market.
$slovakia.Activate()
$slovakia.Range("B4").Select()
